# Rookie Team Season Budget - add two "Game Supplies" line items
# (gobuilda parts / gobuild parts gears) to the Expenses table on Sheet1.
#
# Net effect vs. the original file: two new rows are inserted at row 17
# (pushing everything from the old row 17 onward down by two rows), and the
# new rows are populated with the new line items. Excel auto-adjusts every
# SUM()/reference formula and the merged-cell ranges that are affected by
# the insertion, which is exactly what the target workbook shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two blank rows above the current row 17 ("Extra Battery, Servo
# Adapter, Servo Extension"). This shifts rows 17-50 down to 19-52 and
# updates every formula/merged-cell reference that spans the insertion
# point automatically.
$ws.Rows.Item(17).Resize(2).Insert()

# The new rows should look like row 16 (same borders/number formats).
# Copy formatting column-by-column (not whole-row) so we don't also copy
# formatting into unused columns F:I and bloat the sheet's dimension.
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)

$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B18").PasteSpecial(-4122)

$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("D16").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("E16").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("J16").Copy()
$ws.Range("J17").PasteSpecial(-4122)
$ws.Range("J18").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# Fill in row 18 first and row 17 second, so the shared-string table picks
# up "gobuilda parts" before "gobuild parts gears" (matches upstream order).
$ws.Range("A18").Value = "gobuilda parts"
$ws.Range("B18").Value = 122.8

$ws.Range("A17").Value = "gobuild parts gears"
$ws.Range("B17").Value = 11.64

# Match the workbook's final on-screen selection/scroll position.
$ws.Range("B18").Select()
